$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "42.506.82"
$ws.Range("E2").Value = "  +1.18%  "
Set-TextValue "D3" "2.244.24"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.21%  "
Set-TextValue "D5" "244.80"
$ws.Range("E5").Value = "  -0.96%  "
Set-TextValue "D6" "0.629"
$ws.Range("E6").Value = "  +1.04%  "
Set-TextValue "D7" "75.62"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  +0.05%  "
Set-TextValue "D9" "0.622"
$ws.Range("E9").Value = "  -1.10%  "
Set-TextValue "D10" "43.84"
$ws.Range("E10").Value = "  +7.22%  "
Set-TextValue "D11" "0.0948"
$ws.Range("E11").Value = "  -0.25%  "
Set-TextValue "D12" "7.21"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +0.63%  "
Set-TextValue "D14" "14.58"
$ws.Range("E14").Value = "  -1.50%  "
Set-TextValue "D15" "0.859"
$ws.Range("E15").Value = "  +0.03%  "
Set-TextValue "D16" "2.267.95"
$ws.Range("E16").Value = "  +1.16%  "
Set-TextValue "D17" "42.288.03"
$ws.Range("E17").Value = "  +0.91%  "
Set-TextValue "D18" "0.0000101"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("E19").Value = "  +1.02%  "
Set-TextValue "D20" "72.02"
$ws.Range("E20").Value = "  +0.46%  "
Set-TextValue "D21" "11.12"
$ws.Range("E21").Value = "  +54.60%  "
Set-TextValue "D22" "2.23"
$ws.Range("E22").Value = "  -4.25%  "
Set-TextValue "D23" "231.71"
$ws.Range("E23").Value = "  +0.11%  "
Set-TextValue "D24" "11.76"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +4.67%  "
Set-TextValue "D29" "167.01"
Set-TextValue "D30" "20.69"
$ws.Range("E30").Value = "  +0.66%  "
Set-TextValue "D31" "5.86"
$ws.Range("E31").Value = "  +20.22%  "
Set-TextValue "D32" "0.0816"
$ws.Range("E32").Value = "  -2.18%  "
Set-TextValue "D33" "30.88"
$ws.Range("E33").Value = "  -5.88%  "
Set-TextValue "D34" "0.118"
$ws.Range("E34").Value = "  -2.21%  "
Set-TextValue "D35" "0.125"
$ws.Range("E35").Value = "  -0.33%  "
Set-TextValue "D36" "4.65"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("E37").Value = "  +4.59%  "
Set-TextValue "D38" "13.74"
$ws.Range("E38").Value = "  -1.40%  "
Set-TextValue "D40" "5.76"
$ws.Range("E40").Value = "  -2.70%  "
Set-TextValue "D41" "63.64"
$ws.Range("E41").Value = "  +3.95%  "
$ws.Range("E42").Value = "  -0.28%  "
Set-TextValue "D43" "106.62"
$ws.Range("E43").Value = "  -5.41%  "
Set-TextValue "D44" "8.87"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("E46").Value = "  -0.02%  "
Set-TextValue "D47" "2.42"
$ws.Range("E47").Value = "  +7.12%  "
$ws.Range("E48").Value = "  +0.74%  "
Set-TextValue "D49" "1.18"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +1.29%  "
Set-TextValue "D51" "4.12"
$ws.Range("E51").Value = "  -1.49%  "
